# Add a new "roi_center_region" column to the FieldIlluminationOutput sheet.
# It is inserted right before the existing "processing_application" column
# (column J), shifting processing_application..comment one column to the
# right (J:Q -> K:R) and growing the sheet's used range from A1:Q1 to A1:R1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FieldIlluminationOutput")

# Insert a new blank column at J, pushing existing J:Q columns to K:R.
$ws.Columns("J:J").Insert()

# Populate the header of the newly inserted column.
$ws.Range("J1").Value = "roi_center_region"
